$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("data-input")
$ws2 = $wb.Worksheets.Item("dir_data-input")

# --- Sheet "data-input": column V (path_len:) values reduced by 25 ---
$ws1.Range("V2").Value = 137
$ws1.Range("V3").Value = 161
$ws1.Range("V4").Value = 158
$ws1.Range("V5").Value = 177
$ws1.Range("V6").Value = 127
$ws1.Range("V7").Value = 146
$ws1.Range("V8").Value = 169
$ws1.Range("V9").Value = 134
$ws1.Range("V10").Value = 186
$ws1.Range("V11").Value = 119
$ws1.Range("V12").Value = 170
$ws1.Range("V13").Value = 173

# --- Sheet "dir_data-input": column E (Date_Modified:) tiny precision updates ---
$ws2.Range("E4").Value = 43970.45450548951
$ws2.Range("E7").Value = 43962.69165836073
$ws2.Range("E8").Value = 43962.69165788621
$ws2.Range("E9").Value = 43962.69165755047
$ws2.Range("E10").Value = 43970.45450720222
$ws2.Range("E11").Value = 43962.69159603336
$ws2.Range("E13").Value = 43970.45450803575
